$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the workbook's inline-string cells)
# so numeric-looking strings (e.g. "0.1309") are not coerced into binary
# doubles (which would corrupt trailing zeros / exact decimal text).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '243.59'
Set-TextValue 'D3' '23.82'
Set-TextValue 'D4' '5.259'
Set-TextValue 'D5' '0.05815'
Set-TextValue 'D6' '6.483'
Set-TextValue 'D7' '3.344'
Set-TextValue 'D8' '0.8083'
Set-TextValue 'D9' '0.8779'
Set-TextValue 'D10' '0.1389'
Set-TextValue 'D11' '0.07279'
Set-TextValue 'D13' '0.03056'
Set-TextValue 'D14' '0.09319'
Set-TextValue 'D15' '3.861'
Set-TextValue 'D16' '0.001534'
Set-TextValue 'D17' '0.04694'
Set-TextValue 'D18' '0.0006039'
Set-TextValue 'D19' '0.006154'
Set-TextValue 'D20' '0.001266'
Set-TextValue 'D21' '0.004596'
Set-TextValue 'D22' '0.00008700'
Set-TextValue 'D23' '3.570'
Set-TextValue 'D24' '2.163'
Set-TextValue 'D26' '0.1309'
Set-TextValue 'D28' '0.0002344'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.006330'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1053'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.002621'
$ws.Range('E43').Value = '42CEJICEJI'
Set-TextValue 'D44' '0.008000'
Set-TextValue 'D45' '0.00005541'
Set-TextValue 'D47' '0.5501'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINWorstin24h'
Set-TextValue 'D48' '0.009164'

Write-Host "Applied all 41 cell updates"
